$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently contains a 2-row Table1 (A1:H2 header+data).
# We add a new data row (row 3) for Arulraj Vellingiri's form response,
# growing the table to A1:H3.

$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Copy row 2's formatting down into the new row 3 (keeps the date/time
# number formatting on columns B, C and F consistent with the rest of
# the table).
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122) # xlPasteFormats

# Now fill in the actual response values for the new row.
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 45831.5861226852
$ws.Range("C3").Value = 45831.5866319444
$ws.Range("D3").Value = "avellingiri@hoopp.com"
$ws.Range("E3").Value = "Arulraj Vellingiri"
$ws.Range("G3").Value = "Monday;Thursday;Friday;"
$ws.Range("H3").Value = "16-W625"
